{"js": "// Fix the typo \"incides\" -> \"indices\" in the \"Square brackets...\" bullet,\n// then append a batch of new research-notes bullets directly after it\n// (before the trailing blank paragraphs at the end of the document).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph containing the typo so the script is resilient to\n// any surrounding content.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Square brackets for array incides\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Square brackets for array incides' paragraph\");\n}\n\n// Fix the typo in place, keeping the paragraph/run formatting intact.\ntarget.insertText(\"Square brackets for array indices\", Word.InsertLocation.replace);\nawait context.sync();\n\n// New bullet lines to insert, in order, right after the fixed paragraph.\nconst newLines = [\n  \"Apps can create interactive visualization of models\",\n  \"Apps allow for adjustments to parameters, displaying resulting changes\",\n  \"Who used similar models to track ebola outbreak\",\n  \"Stochastic processes should be used to account for randomness and probability\",\n  \"Poisson distribution = describes probability of certain number of events occurring in fixed interval of time and space\",\n  \"Can be used to model spread of epidemic\",\n  \"Mean of poisson = average number of events in given interval - determines shape of distribution curve\",\n  \"Poisson distribution used to model change in number of susceptible, infected and dead individual over given time step in context of epidemic model\",\n  \"Covid 19 provided useful data for future analysis\",\n  \"Uk requires death certificates therefore covid 19 mortalitiy rates are reliable\",\n  \"Sir equations used to estimate covid 19 infections during pandemic\",\n  \"Rate of change of deaths  = proportional to number of infected individuals\",\n  \"Linear interpolation necessary to evaluate models at specific times\"\n];\n\nlet anchor = target;\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Fix the typo \"incides\" -> \"indices\" in the \"Square brackets...\" bullet,\n# then append a batch of new research-notes bullets directly after it\n# (before the trailing blank paragraphs at the end of the document).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph with the typo so the script is resilient to any\n# surrounding content, then fix it by writing directly into its Range\n# (keeps the existing run/paragraph formatting intact).\n$targetIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Square brackets for array incides*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.Text = \"Square brackets for array indices\"\n\n$newLines = @(\n    \"Apps can create interactive visualization of models\",\n    \"Apps allow for adjustments to parameters, displaying resulting changes\",\n    \"Who used similar models to track ebola outbreak\",\n    \"Stochastic processes should be used to account for randomness and probability\",\n    \"Poisson distribution = describes probability of certain number of events occurring in fixed interval of time and space\",\n    \"Can be used to model spread of epidemic\",\n    \"Mean of poisson = average number of events in given interval - determines shape of distribution curve\",\n    \"Poisson distribution used to model change in number of susceptible, infected and dead individual over given time step in context of epidemic model\",\n    \"Covid 19 provided useful data for future analysis\",\n    \"Uk requires death certificates therefore covid 19 mortalitiy rates are reliable\",\n    \"Sir equations used to estimate covid 19 infections during pandemic\",\n    \"Rate of change of deaths  = proportional to number of infected individuals\",\n    \"Linear interpolation necessary to evaluate models at specific times\"\n)\n\n$idx = $targetIndex\nforeach ($line in $newLines) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter()\n    $idx = $idx + 1\n    $newP = $d.Paragraphs.Item($idx)\n    $newP.Range.Text = $line\n}\n"}
